# Import additional quiz questions/answers into the sheet.
# Question 2's answer options are regenerated with a "2_" prefix, and two
# new questions (3 and 4) with their answers are appended below the
# existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Question 2 block (rows 6-10): refresh answer text ---------------------
$ws.Range("B6").Value = "2_Не правильный ответ 4"
$ws.Range("B7").Value = "2_Не правильный ответ 1"
$ws.Range("B8").Value = "2_Не правильный ответ 2"
$ws.Range("B9").Value = "2_Правильный ответ"
$ws.Range("C9").Value = 1
$ws.Range("B10").Value = "2_Не правильный ответ 3"
$ws.Range("C10").ClearContents()

# --- Question 3 block (rows 11-12) ------------------------------------------
$ws.Range("A11").Value = "Пробный вопрос 3"
$ws.Range("B11").Value = "3_Правильный ответ"
$ws.Range("C11").Value = 1
$ws.Range("B12").Value = "3_Не правильный ответ 1"

# --- Question 4 block (rows 13-15) ------------------------------------------
$ws.Range("A13").Value = "Пробный вопрос 4"
$ws.Range("B13").Value = "4_Не правильный ответ 1"
$ws.Range("A14").Value = " "
$ws.Range("B14").Value = "4_Правильный ответ"
$ws.Range("C14").Value = 1
$ws.Range("B15").Value = "4_Не правильный ответ 2"

# Column B is wider now to fit the longer imported answer text.
$ws.Columns.Item(2).ColumnWidth = 30.59

# Leave the selection where the import left off.
$ws.Range("B11").Select()
